# CIERRE DEL 29 OCT 21
# Record payments received against existing credit remisiones on the
# "REMISIONES OCTUBRE  2021     " sheet (rows 10-14), then leave the
# selection on D15 as in the source workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)
$ws.Activate()

# Row 10 - Remision #44 (10460) paid 17-Oct-21
$ws.Range("F10").Value = 44486
$ws.Range("G10").Value = 10460

# Row 11 - Remision #45 (7637) paid 17-Oct-21
$ws.Range("F11").Value = 44486
$ws.Range("G11").Value = 7637

# Row 12 - new remision dated 17-Oct-21, GUSTAVO, 3631, paid 19-Oct-21
$ws.Range("A12").Value = 44486
$ws.Range("D12").Value = "GUSTAVO"
$ws.Range("E12").Value = 3631
$ws.Range("F12").Value = 44488
$ws.Range("G12").Value = 3631

# Row 13 - new remision dated 19-Oct-21, GUSTAVO, 4277, paid 21-Oct-21
$ws.Range("A13").Value = 44488
$ws.Range("D13").Value = "GUSTAVO"
$ws.Range("E13").Value = 4277
$ws.Range("F13").Value = 44490
$ws.Range("G13").Value = 4277

# Row 14 - new remision dated 22-Oct-21, GUSTAVO, 5290 (still unpaid)
$ws.Range("A14").Value = 44491
$ws.Range("D14").Value = "GUSTAVO"
$ws.Range("E14").Value = 5290

# Restore the selection recorded in the workbook after this edit
$ws.Range("D15").Select()
